# Update the check_ssl_cert stats workbook after the latest release:
# append a new data row (2022-11-18 / serial 44840) to the "Data" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Grow the "Data" table by one row; this also extends the worksheet
# dimension and the table/autoFilter reference (A1:AH47 -> A1:AH48).
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the formatting (number formats / styles) of the last existing
# row down into the freshly added row 48 before filling in values.
$ws.Range("A47:AH47").Copy() | Out-Null
$ws.Range("A48:AH48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false | Out-Null

# Plain data values for row 48.
$ws.Range("A48").Value2 = 44840
$ws.Range("B48").Value2 = 106
$ws.Range("C48").Value2 = 255
$ws.Range("D48").Value2 = 207
$ws.Range("E48").Value2 = 5260
$ws.Range("G48").Value2 = 6107
$ws.Range("H48").Value2 = 1875
$ws.Range("I48").Value2 = 291
$ws.Range("J48").Value2 = 425
$ws.Range("K48").Value2 = 98
$ws.Range("L48").Value2 = 57
$ws.Range("O48").Value2 = 1890
$ws.Range("P48").Value2 = 4001
$ws.Range("Q48").Value2 = 64679
$ws.Range("R48").Value2 = 44879
$ws.Range("S48").Value2 = 0
$ws.Range("T48").Value2 = 0
$ws.Range("U48").Value2 = 238
$ws.Range("W48").Value2 = 0
$ws.Range("X48").Value2 = 162
$ws.Range("Z48").Value2 = 151
$ws.Range("AA48").Value2 = 161
$ws.Range("AB48").Value2 = 4
$ws.Range("AC48").Value2 = 0
$ws.Range("AD48").Value2 = 345
$ws.Range("AE48").Value2 = 837
$ws.Range("AF48").Value2 = 9

# Calculated ("formula") columns of the table, same structured
# references as the rest of the table.
$ws.Range("F48").Formula = "=Data[[#This Row],[LoC]]-E47"
$ws.Range("M48").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("N48").Formula = "=Data[[#This Row],[Total]]-M47"
$ws.Range("V48").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("Y48").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AH48").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# Leave [GH runs] (AG48) empty, same as the source row template.

# Reflect where the user ended up after typing in the new row.
$ws.Range("AI48").Select() | Out-Null
